# Generate Report for Handback
# The 875e3f4d-... and ce59ca7d-... files have now been handed back
# (target files generated) for both the zh-cn and de-de locales.
# Update the Overview sheet and the per-locale detail sheets to reflect
# the new "Handed back: in sync with en-US" status, populate the
# "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns, and add hyperlinks for the newly-produced target
# files.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: rows 4 (875e3f4d...) and 5 (ce59ca7d...) - the zh-cn
# (E) and de-de (F) status columns move from "Ready for handoff" to
# "Handed back: in sync with en-US".
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E4").Value = $handedBack
$wsOverview.Range("F4").Value = $handedBack
$wsOverview.Range("E5").Value = $handedBack
$wsOverview.Range("F5").Value = $handedBack

# ---------------------------------------------------------------------
# zh-cn sheet: rows 4 (875e3f4d...) and 5 (ce59ca7d...)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Row 4 - 875e3f4d-9958-4150-968d-4fe972d513ab
$wsZhCn.Range("C4").Value = $handedBack
$wsZhCn.Range("J4").Value = $wsZhCn.Range("G4").Value2
$wsZhCn.Range("K4").Value = "2016-08-31 12:35:46"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/8248db31a73fb01048f4b927b751a97f3b12c62d/e2e/875e3f4d-9958-4150-968d-4fe972d513ab.md",
    "",
    "",
    "875e3f4d-9958-4150-968d-4fe972d513ab.md"
)

# Row 5 - ce59ca7d-2df5-4a7d-9f37-c85d1b0d44f6
$wsZhCn.Range("C5").Value = $handedBack
$wsZhCn.Range("J5").Value = $wsZhCn.Range("G5").Value2
$wsZhCn.Range("K5").Value = "2016-08-31 12:35:46"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I5"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/8248db31a73fb01048f4b927b751a97f3b12c62d/e2e/ce59ca7d-2df5-4a7d-9f37-c85d1b0d44f6.md",
    "",
    "",
    "ce59ca7d-2df5-4a7d-9f37-c85d1b0d44f6.md"
)

# ---------------------------------------------------------------------
# de-de sheet: rows 4 (875e3f4d...) and 5 (ce59ca7d...)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 4 - 875e3f4d-9958-4150-968d-4fe972d513ab
$wsDeDe.Range("C4").Value = $handedBack
$wsDeDe.Range("J4").Value = $wsDeDe.Range("G4").Value2
$wsDeDe.Range("K4").Value = "2016-08-31 12:35:54"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/828591943ea3f3fda32733380aba6a14a1043ce2/e2e/875e3f4d-9958-4150-968d-4fe972d513ab.md",
    "",
    "",
    "875e3f4d-9958-4150-968d-4fe972d513ab.md"
)

# Row 5 - ce59ca7d-2df5-4a7d-9f37-c85d1b0d44f6
$wsDeDe.Range("C5").Value = $handedBack
$wsDeDe.Range("J5").Value = $wsDeDe.Range("G5").Value2
$wsDeDe.Range("K5").Value = "2016-08-31 12:35:54"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I5"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/828591943ea3f3fda32733380aba6a14a1043ce2/e2e/ce59ca7d-2df5-4a7d-9f37-c85d1b0d44f6.md",
    "",
    "",
    "ce59ca7d-2df5-4a7d-9f37-c85d1b0d44f6.md"
)
